$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128; this shifts existing rows 128-207 down
# to 129-208 (matching the diff: every record from old row 128 onward moves
# down by one row, and a brand-new record is inserted at row 128).
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new record's data.
$ws.Cells.Item(128, 1).Value = 10
$ws.Cells.Item(128, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(128, 3).Value = "La Araucanía"
$ws.Cells.Item(128, 4).Value = 45068
$ws.Cells.Item(128, 5).Value = 9
$ws.Cells.Item(128, 6).Value = 100112031
$ws.Cells.Item(128, 7).Value = "Poroto verde"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 50
$ws.Cells.Item(128, 11).Value = 20000
$ws.Cells.Item(128, 12).Value = 20000
$ws.Cells.Item(128, 13).Value = 20000
$ws.Cells.Item(128, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(128, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(128, 16).Value = 800
$ws.Cells.Item(128, 17).Value = 25
$ws.Cells.Item(128, 18).Value = "Hortaliza"
